# ---------------------------------------------------------------------------
# LC_columns.xlsx update
#   - fix two "X.0 µm" -> "X µm" typos in particle-size descriptions
#   - rename the Bruker MOSAIC "Class" label and rewrite its descriptions
#   - turn the D/E (Length / ID) columns into text-formatted columns
#   - add 7 new rows describing Waters Acquity, Inertsil trap and
#     IonOptics Aurora columns
#   - tidy up the sheet view (zoom, selection) and column widths
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Small text corrections on existing rows
# ---------------------------------------------------------------------------
$ws.Range("I2").Value = "EasySpray PepMap RSLC column (2 µm C18-coated particles, 50 cm * 75 µm ID, ThermoFisher Scientific P/N ES903)"
$ws.Range("I5").Value = "Acclaim PepMap trap column (5 µm C18-coated particles, 0.5 cm * 300 µm ID, ThermoFisher Scientific P/N 160454)"

# ---------------------------------------------------------------------------
# 2) Bruker monolithic / MOSAIC rows (10 & 11): rename Class + rewrite text
# ---------------------------------------------------------------------------
$ws.Range("B10").Value = "Bruker C18 MOSAIC column"
$ws.Range("I10").Value = "experimental 15 cm monolithic capillary column (Bruker Daltonics) column (C18-coated particles, 15 cm)"
$ws.Range("D10").Value = "15"

$ws.Range("B11").Value = "Bruker C18 MOSAIC column"
$ws.Range("D11").Value = "15"
$ws.Range("E11").Value = "75"
$ws.Range("I11").Value = "monolithic capillary MOSAIC column (C18-coated particles, 15 cm * 75 µm ID, Bruker P/N 1908374)"
$ws.Range("J11").Value = "1908374"

# ---------------------------------------------------------------------------
# 3) D & E columns (Length.(cm) / ID.(µm)) become text for all existing data
#    rows (2-14): apply a Text number format, matching the values already
#    present (numbers stored as text strings).
# ---------------------------------------------------------------------------
$ws.Range("D2:E14").NumberFormat = "@"

$ws.Range("D2").Value = "50"
$ws.Range("E2").Value = "75"
$ws.Range("D3").Value = "200"
$ws.Range("D4").Value = "50"
$ws.Range("D5").Value = "0.5"
$ws.Range("E5").Value = "300"
$ws.Range("D7").Value = "25"
$ws.Range("E7").Value = "150"
$ws.Range("D8").Value = "10"
$ws.Range("E8").Value = "150"
$ws.Range("D9").Value = "10"
$ws.Range("E9").Value = "75"

# ---------------------------------------------------------------------------
# 4) New rows 12-14 (Waters Acquity BEH columns + Inertsil ODS-4 trap)
# ---------------------------------------------------------------------------
$ws.Range("A12").Value = "Acquity UPLC Peptide BEH C18 column, 15 cm long"
$ws.Range("B12").Value = "Acquity Peptide BEH C18 column"
$ws.Range("C12").Value = "Waters"
$ws.Range("D12").Value = "15"
$ws.Range("E12").Value = "2100"
$ws.Range("F12").Value = "1.7"
$ws.Range("G12").Value = "C18"
$ws.Range("H12").Value = "packed bed"
$ws.Range("I12").Value = "Acquity UPLC Peptide BEH C18 column (1.7 µm C18-coated particles, 15 cm * 2100 µm ID, Waters P/N 186002353)"
$ws.Range("J12").Value = "186002353"
$ws.Range("K12").Value = "Fractionation"

$ws.Range("A13").Value = "Acquity UPLC Peptide BEH C18 column, 30 cm long"
$ws.Range("B13").Value = "Acquity Peptide BEH C18 column"
$ws.Range("C13").Value = "Waters"
$ws.Range("D13").Value = "30"
$ws.Range("E13").Value = "2100"
$ws.Range("F13").Value = "1.7"
$ws.Range("G13").Value = "C18"
$ws.Range("H13").Value = "packed bed"
$ws.Range("I13").Value = "Acquity UPLC Peptide BEH C18 column (1.7 µm C18-coated particles, 30 cm * 2100 µm ID, Waters P/N 186005792)"
$ws.Range("J13").Value = "186005792"
$ws.Range("K13").Value = "Fractionation"

$ws.Range("A14").Value = "experimental Inertsil ODS-4"
$ws.Range("B14").Value = "Intersil ODS4 trap column"
$ws.Range("C14").Value = "Bruker"
$ws.Range("D14").Value = "3"
$ws.Range("G14").Value = "C18"
$ws.Range("H14").Value = "packed bed"
$ws.Range("I14").Value = "experimental Inertsil ODS-4 trap column (C18-coated particles, 3 cm)"
$ws.Range("K14").Value = "Trap"

$ws.Range("D14:E14").NumberFormat = "@"

# ---------------------------------------------------------------------------
# 5) New rows 15-18 (IonOptics Aurora columns) - D/E hold real numbers, with
#    a Text number format applied (matches how the author padded these
#    cells: the numbers display as plain digits).
# ---------------------------------------------------------------------------
$ws.Range("A15").Value = "Aurora Ultimate CSI 25 cm column"
$ws.Range("B15").Value = "Aurora column"
$ws.Range("C15").Value = "IonOptics"
$ws.Range("D15").Value = 25
$ws.Range("E15").Value = 75
$ws.Range("F15").Value = "1.7"
$ws.Range("G15").Value = "C18"
$ws.Range("H15").Value = "packed bed"
$ws.Range("I15").Value = "Aurora Ultimate CSI column (C18-coated particles, 25 cm * 75 C18 UHPLC column"
$ws.Range("J15").Value = "Aurora Ultimate CSI 25 cm column"
$ws.Range("K15").Value = "Analytical"

$ws.Range("A16").Value = "Aurora Ultimate CSI 60 cm column"
$ws.Range("B16").Value = "Aurora column"
$ws.Range("C16").Value = "IonOptics"
$ws.Range("D16").Value = 60
$ws.Range("E16").Value = 75
$ws.Range("F16").Value = "1.7"
$ws.Range("G16").Value = "C18"
$ws.Range("H16").Value = "packed bed"
$ws.Range("I16").Value = "Aurora Ultimate CSI column (C18-coated particles, 60 cm * 75 C18 UHPLC column"
$ws.Range("J16").Value = "Aurora Ultimate CSI 60 cm column"
$ws.Range("K16").Value = "Analytical"

$ws.Range("A17").Value = "Aurora Ultimate XT 25 cm column"
$ws.Range("B17").Value = "Aurora column"
$ws.Range("C17").Value = "IonOptics"
$ws.Range("D17").Value = 25
$ws.Range("E17").Value = 75
$ws.Range("F17").Value = "1.7"
$ws.Range("G17").Value = "C18"
$ws.Range("H17").Value = "packed bed"
$ws.Range("I17").Value = "Aurora Ultimate XT column (C18-coated particles, 25 cm * 75 C18 UHPLC column"
$ws.Range("J17").Value = "Aurora Ultimate XT 25 cm column"
$ws.Range("K17").Value = "Analytical"

$ws.Range("A18").Value = "Aurora Ultimate XT 60 cm column"
$ws.Range("B18").Value = "Aurora column"
$ws.Range("C18").Value = "IonOptics"
$ws.Range("D18").Value = 60
$ws.Range("E18").Value = 75
$ws.Range("F18").Value = "1.7"
$ws.Range("G18").Value = "C18"
$ws.Range("H18").Value = "packed bed"
$ws.Range("I18").Value = "Aurora Ultimate XT column (C18-coated particles, 60 cm * 75 C18 UHPLC column"
$ws.Range("J18").Value = "Aurora Ultimate XT 60 cm column"
$ws.Range("K18").Value = "Analytical"

$ws.Range("D15:E18").NumberFormat = "@"

# ---------------------------------------------------------------------------
# 6) P/N column (J) keeps a Text format for the new rows that reuse it (10,14)
# ---------------------------------------------------------------------------
$ws.Range("J10").NumberFormat = "@"
$ws.Range("J14").NumberFormat = "@"

# ---------------------------------------------------------------------------
# 7) Column widths (best-fit on the new/changed columns)
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 61.28515625
$ws.Columns.Item(2).ColumnWidth = 35.28515625
$ws.Columns.Item(3).ColumnWidth = 22.28515625
$ws.Columns.Item(4).ColumnWidth = 11.5703125
$ws.Columns.Item(5).ColumnWidth = 7.7109375
$ws.Columns.Item(6).ColumnWidth = 17.7109375
$ws.Columns.Item(7).ColumnWidth = 8.5703125
$ws.Columns.Item(10).ColumnWidth = 33.5703125

# ---------------------------------------------------------------------------
# 8) Sheet view: zoom to 70%, move the selection, drop the frozen top-left
# ---------------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 70
$ws.Range("A21").Select() | Out-Null
